# Auto-generated Excel COM-interop script applying the Durandal_Profits
# workbook update (H/I/J/K/L/M/N recompute across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
# Row 64
$ws.Range("H64").Value = 3741.5
$ws.Range("I64").Value = 3982
$ws.Range("J64").Value = 3693.4
$ws.Range("K64").Value = 3982
$ws.Range("L64").Value = 3693.4
$ws.Range("M64").Value = -3734
$ws.Range("N64").Value = -4189.4
# Row 67
$ws.Range("H67").Value = 3741.5
$ws.Range("I67").Value = 3982
$ws.Range("J67").Value = 3693.4
$ws.Range("K67").Value = 3982
$ws.Range("L67").Value = 3693.4
$ws.Range("M67").Value = -3124
$ws.Range("N67").Value = -5409.4
# Row 74
$ws.Range("H74").Value = 3628.111
$ws.Range("I74").Value = 2986.125
$ws.Range("J74").Value = 4141.7
$ws.Range("K74").Value = 2986.125
$ws.Range("L74").Value = 4141.7
$ws.Range("M74").Value = -2050.125
$ws.Range("N74").Value = -6013.7
# Row 76
$ws.Range("H76").Value = 2472686.8
$ws.Range("I76").Value = 2648910.5
$ws.Range("K76").Value = 2648910.5
$ws.Range("M76").Value = -2648595.5
# Row 77
$ws.Range("H77").Value = 3628.111
$ws.Range("I77").Value = 2986.125
$ws.Range("J77").Value = 4141.7
$ws.Range("K77").Value = 14930.625
$ws.Range("L77").Value = 20708.5
$ws.Range("M77").Value = -10250.625
$ws.Range("N77").Value = -30068.5
# Row 79
$ws.Range("H79").Value = 2472686.8
$ws.Range("I79").Value = 2648910.5
$ws.Range("K79").Value = 2648910.5
$ws.Range("M79").Value = -2647818.5
# Row 135
$ws.Range("H135").Value = 2732.1333
$ws.Range("I135").Value = 1406.0741
$ws.Range("J135").Value = 14666.667
$ws.Range("K135").Value = 12654.6669
$ws.Range("L135").Value = 132000.003
$ws.Range("M135").Value = -10119.6669
$ws.Range("N135").Value = -137070.003
# Row 136
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()
# Row 137
$ws.Range("H137").Value = 1136.9062
$ws.Range("I137").Value = 834.15
$ws.Range("J137").Value = 1641.5
$ws.Range("K137").Value = 2502.45
$ws.Range("L137").Value = 4924.5
$ws.Range("M137").Value = 47.55000000000018
$ws.Range("N137").Value = -10024.5
# Row 139
$ws.Range("H139").Value = 70111.42999999999
$ws.Range("J139").Value = 70111.42999999999
$ws.Range("L139").Value = 70111.42999999999
$ws.Range("N139").Value = -80391.42999999999
# Row 140
$ws.Range("H140").Value = 68455
$ws.Range("J140").Value = 89221.42999999999
$ws.Range("L140").Value = 89221.42999999999
$ws.Range("N140").Value = -99581.42999999999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 397268.88
$ws.Range("I32").Value = 5442.742
$ws.Range("J32").Value = 4091629.8
$ws.Range("K32").Value = 5442.742
$ws.Range("L32").Value = 4091629.8
$ws.Range("M32").Value = -5155.742
$ws.Range("N32").Value = -4092203.8
# Row 61
$ws.Range("H61").Value = 2517.7576
$ws.Range("I61").Value = 2724.348
$ws.Range("J61").Value = 2042.6
$ws.Range("K61").Value = 2724.348
$ws.Range("L61").Value = 2042.6
$ws.Range("M61").Value = -2512.348
$ws.Range("N61").Value = -2466.6
# Row 74
$ws.Range("H74").Value = 1020.26666
$ws.Range("I74").Value = 647.4286
$ws.Range("J74").Value = 1346.5
$ws.Range("K74").Value = 647.4286
$ws.Range("L74").Value = 1346.5
$ws.Range("M74").Value = 226.5714
$ws.Range("N74").Value = -3094.5
# Row 77
$ws.Range("H77").Value = 1020.26666
$ws.Range("I77").Value = 647.4286
$ws.Range("J77").Value = 1346.5
$ws.Range("K77").Value = 3237.143
$ws.Range("L77").Value = 6732.5
$ws.Range("M77").Value = 1130.857
$ws.Range("N77").Value = -15468.5
# Row 95
$ws.Range("H95").Value = 22000
$ws.Range("J95").Value = 22000
$ws.Range("L95").Value = 22000
$ws.Range("N95").Value = -27492
# Row 135
$ws.Range("H135").Value = 55952.668
$ws.Range("J135").Value = 55952.668
$ws.Range("L135").Value = 55952.668
$ws.Range("N135").Value = -66092.66800000001
# Row 136
$ws.Range("H136").Value = 2517.7576
$ws.Range("I136").Value = 2724.348
$ws.Range("J136").Value = 2042.6
$ws.Range("K136").Value = 8173.044
$ws.Range("L136").Value = 6127.799999999999
$ws.Range("M136").Value = -5623.044
$ws.Range("N136").Value = -11227.8
# Row 138
$ws.Range("H138").Value = 60712.5
$ws.Range("J138").Value = 60712.5
$ws.Range("L138").Value = 60712.5
$ws.Range("N138").Value = -70992.5
# Row 141
$ws.Range("H141").Value = 59800
$ws.Range("J141").Value = 62707.69
$ws.Range("L141").Value = 62707.69
$ws.Range("N141").Value = -73067.69

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 14398.042
$ws.Range("I86").Value = 20621.625
$ws.Range("J86").Value = 1950.875
$ws.Range("K86").Value = 20621.625
$ws.Range("L86").Value = 1950.875
$ws.Range("M86").Value = -19498.625
$ws.Range("N86").Value = -4196.875
# Row 89
$ws.Range("H89").Value = 14398.042
$ws.Range("I89").Value = 20621.625
$ws.Range("J89").Value = 1950.875
$ws.Range("K89").Value = 103108.125
$ws.Range("L89").Value = 9754.375
$ws.Range("M89").Value = -97492.125
$ws.Range("N89").Value = -20986.375
# Row 94
$ws.Range("H94").Value = 964.0571
$ws.Range("I94").Value = 877.29034
$ws.Range("J94").Value = 1636.5
$ws.Range("K94").Value = 877.29034
$ws.Range("L94").Value = 1636.5
$ws.Range("M94").Value = -426.29034
$ws.Range("N94").Value = -2538.5
# Row 105
$ws.Range("H105").Value = 1942.3529
$ws.Range("I105").Value = 1768.3334
$ws.Range("J105").Value = 2360
$ws.Range("K105").Value = 1768.3334
$ws.Range("L105").Value = 2360
$ws.Range("M105").Value = -21.33339999999998
$ws.Range("N105").Value = -5854
# Row 107
$ws.Range("H107").Value = 4089.919
$ws.Range("I107").Value = 4434.5483
$ws.Range("K107").Value = 4434.5483
$ws.Range("M107").Value = -2514.5483
# Row 134
$ws.Range("H134").Value = 4897.5127
$ws.Range("I134").Value = 823.2759
$ws.Range("J134").Value = 16712.8
$ws.Range("K134").Value = 2469.8277
$ws.Range("L134").Value = 50138.39999999999
$ws.Range("M134").Value = 65.17230000000018
$ws.Range("N134").Value = -55208.39999999999
# Row 137
$ws.Range("H137").Value = 72152.92999999999
$ws.Range("J137").Value = 72152.92999999999
$ws.Range("L137").Value = 72152.92999999999
$ws.Range("N137").Value = -82352.92999999999
# Row 138
$ws.Range("H138").Value = 66600
$ws.Range("J138").Value = 66600
$ws.Range("L138").Value = 66600
$ws.Range("N138").Value = -76880
# Row 140
$ws.Range("H140").Value = 73266.664
$ws.Range("J140").Value = 73266.664
$ws.Range("L140").Value = 73266.664
$ws.Range("N140").Value = -83626.664

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 4413.0815
$ws.Range("I31").Value = 2907.861
$ws.Range("J31").Value = 5496.84
$ws.Range("K31").Value = 2907.861
$ws.Range("L31").Value = 5496.84
$ws.Range("M31").Value = -2612.861
$ws.Range("N31").Value = -6086.84
# Row 34
$ws.Range("H34").Value = 4413.0815
$ws.Range("I34").Value = 2907.861
$ws.Range("J34").Value = 5496.84
$ws.Range("K34").Value = 2907.861
$ws.Range("L34").Value = 5496.84
$ws.Range("M34").Value = -2705.861
$ws.Range("N34").Value = -5900.84
# Row 138
$ws.Range("H138").Value = 48200
$ws.Range("J138").Value = 48200
$ws.Range("L138").Value = 48200
$ws.Range("N138").Value = -58480

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 5435655
$ws.Range("I131").Value = 1220.3077
$ws.Range("J131").Value = 6329929
$ws.Range("K131").Value = 3660.9231
$ws.Range("L131").Value = 18989787
$ws.Range("M131").Value = 1379.0769
$ws.Range("N131").Value = -18999867

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 14069543
$ws.Range("I70").Value = 24463256
$ws.Range("J70").Value = 7459.8823
$ws.Range("K70").Value = 24463256
$ws.Range("L70").Value = 7459.8823
$ws.Range("M70").Value = -24462986
$ws.Range("N70").Value = -7999.8823
# Row 73
$ws.Range("H73").Value = 14069543
$ws.Range("I73").Value = 24463256
$ws.Range("J73").Value = 7459.8823
$ws.Range("K73").Value = 24463256
$ws.Range("L73").Value = 7459.8823
$ws.Range("M73").Value = -24462320
$ws.Range("N73").Value = -9331.882300000001
# Row 80
$ws.Range("H80").Value = 64528.11
$ws.Range("I80").Value = 103537
$ws.Range("J80").Value = 3228.4285
$ws.Range("K80").Value = 103537
$ws.Range("L80").Value = 3228.4285
$ws.Range("M80").Value = -102539
$ws.Range("N80").Value = -5224.4285
# Row 83
$ws.Range("H83").Value = 64528.11
$ws.Range("I83").Value = 103537
$ws.Range("J83").Value = 3228.4285
$ws.Range("K83").Value = 517685
$ws.Range("L83").Value = 16142.1425
$ws.Range("M83").Value = -512693
$ws.Range("N83").Value = -26126.1425
# Row 138
$ws.Range("H138").Value = 68366.664
$ws.Range("J138").Value = 68366.664
$ws.Range("L138").Value = 68366.664
$ws.Range("N138").Value = -78646.664
# Row 140
$ws.Range("H140").Value = 89923
$ws.Range("J140").Value = 89923
$ws.Range("L140").Value = 89923
$ws.Range("N140").Value = -100283

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 136
$ws.Range("H136").Value = 2974.6072
$ws.Range("I136").Value = 3400.7273
$ws.Range("J136").Value = 2870.4443
$ws.Range("K136").Value = 10202.1819
$ws.Range("L136").Value = 8611.332900000001
$ws.Range("M136").Value = -7652.1819
$ws.Range("N136").Value = -13711.3329
# Row 138
$ws.Range("H138").Value = 49849
$ws.Range("J138").Value = 49849
$ws.Range("L138").Value = 49849
$ws.Range("N138").Value = -60129
# Row 141
$ws.Range("H141").Value = 65785
$ws.Range("J141").Value = 65785
$ws.Range("L141").Value = 65785
$ws.Range("N141").Value = -76145

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 726.64514
$ws.Range("I107").Value = 781.619
$ws.Range("J107").Value = 611.2
$ws.Range("K107").Value = 2344.857
$ws.Range("L107").Value = 1833.6
$ws.Range("M107").Value = -424.857
$ws.Range("N107").Value = -5673.6
# Row 136
$ws.Range("H136").Value = 759.8919
$ws.Range("I136").Value = 737.2
$ws.Range("K136").Value = 2211.6
$ws.Range("M136").Value = 338.3999999999996
# Row 138
$ws.Range("H138").Value = 73600
$ws.Range("J138").Value = 73600
$ws.Range("L138").Value = 73600
$ws.Range("N138").Value = -83880
# Row 139
$ws.Range("H139").Value = 61843
$ws.Range("J139").Value = 61843
$ws.Range("L139").Value = 61843
$ws.Range("N139").Value = -72123
# Row 141
$ws.Range("H141").Value = 75943
$ws.Range("J141").Value = 75943
$ws.Range("L141").Value = 75943
$ws.Range("N141").Value = -86303

